# Fixed serializer, glob order lookup
# Rebuild the order-line rows: correct SKU/name alignment, fix quantities
# for Fiji Water and Guayaki Enlighten Mint, and append the two missing
# Guayaki Yerba Mate lines (Revel Berry, Bluephoria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-7 (SKU, Name, Quantity, Cost Per, Total Cost)
$rows = @(
    @("77802", "Ithaca Soda - Ginger Beer", "1", "23.95", "23.95"),
    @("77801", "Ithaca Soda - Root Beer", "1", "23.95", "23.95"),
    @("75112", "Fiji Water 1L", "3", "16.50", "49.50"),
    @("77220", "Guayaki Yerba Mate - Enlighten Mint", "2", "28.00", "56.00"),
    @("77221", "Guayaki Yerba Mate - Revel Berry", "2", "28.00", "56.00"),
    @("77222", "Guayaki Yerba Mate - Bluephoria", "2", "28.00", "56.00")
)

# Keep SKU/Quantity/Cost/Total columns stored as text, matching how this
# sheet already represents every value (inline strings, not numbers).
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("C2:E7").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
